# Apply the LOM3207 edit.
#
# Net effect vs. the original workbook:
#  - Row 13 (the old "Docentes responsáveis:" value row, which only had
#    B13/C13 populated) is removed, shifting rows 14-24 up to rows 13-23.
#  - The long "Objetivos:" paragraph in B10/C10 is overwritten with the
#    teacher name string.
#  - A handful of the long descriptive paragraphs that land in the shifted
#    rows get overwritten with shorter (reused) values instead of keeping
#    their original long text.
#
# Strategy: capture the handful of reused text values up front (before any
# mutation), then write them into their destinations, and finally delete
# row 13 so everything slides up into place. Writing via Range.Copy (rather
# than Range.Value = "...") for values that already exist elsewhere in the
# sheet keeps them as plain shared-string text and avoids Excel's automatic
# "this looks like a date" reinterpretation/style churn that a fresh
# .Value assignment of "01/01/2012" would trigger.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the source cells we'll need, before any writes happen.
$teacherName = $ws.Range("B13").Value()        # "5982760 - Carlos Alberto Baldan"
$activationDate = $ws.Range("B8").Value()      # "01/01/2012"
$metodoText = $ws.Range("B19").Value()         # "Aulas expositivas e práticas..."
$criterioText = $ws.Range("B20").Value()       # "Média de duas provas..."
$normaText = $ws.Range("B21").Value()          # "Aplicação de uma prova..."

# 1) Overwrite the long "Objetivos:" description with the teacher's name.
$ws.Range("B13").Copy($ws.Range("B10"))
$ws.Range("C13").Copy($ws.Range("C10"))

# 2) Fix up the cells that will land in rows 13, 15, 18, 19, 20, 21 after
#    row 13 is deleted -- write into their pre-delete (current) row numbers.

# -> after-row 13 (currently row 14): brand-new text "Semestral".
$ws.Range("B14").Value = "Semestral"
$ws.Range("C14").Value = "Semestral"

# -> after-row 15 (currently row 16): reuse the activation date text.
$ws.Range("B8").Copy($ws.Range("B16"))
$ws.Range("C8").Copy($ws.Range("C16"))

# -> after-row 18 (currently row 19): reuse the teacher name text again.
$ws.Range("B19").Value = $teacherName
$ws.Range("C19").Value = $teacherName

# -> after-row 19 (currently row 20): reuse the "Método" text.
$ws.Range("B20").Value = $metodoText
$ws.Range("C20").Value = $metodoText

# -> after-row 20 (currently row 21): reuse the "Critério" text.
$ws.Range("B21").Value = $criterioText
$ws.Range("C21").Value = $criterioText

# -> after-row 21 (currently row 22): reuse the "Norma de recuperação" text.
$ws.Range("B22").Value = $normaText
$ws.Range("C22").Value = $normaText

# 3) Finally delete row 13, shifting rows 14-24 up to rows 13-23 and
#    dropping the sheet's used range to A1:C23.
$ws.Rows(13).Delete()

Write-Output "done"
